# first sweep cleaning data columns to conform to specs--done by chase
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Normalize free-text category labels to underscore-separated tokens.
$used.Replace("Environmental Perturbation", "Environmental_Perturbation") | Out-Null
$used.Replace("KN99 alpha", "KN99_alpha") | Out-Null
$used.Replace("Time Course", "Timecourse") | Out-Null

# Move the view/selection to where the editor last left off.
$ws.Range("F36").Select()
